# Applies the diff described for UC007 worksheet:
#  - Bump the "Version:" value in D2 from 0.1 to 1.0
#  - Rotate the second step (row 20 / row 28 / row 36) text among TC2, TC3, TC4
#    so that:
#      TC2's 2nd step becomes what used to be TC4's 2nd step
#      TC3's 2nd step becomes what used to be TC2's 2nd step
#      TC4's 2nd step becomes what used to be TC3's 2nd step

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump version number (leading apostrophe keeps "1.0" stored as text,
#    matching the original "0.1" text value instead of being coerced to
#    the number 1)
$ws.Range("D2").Value = "'1.0"

# 2. Capture current (pre-edit) values before overwriting anything
$tc2_B20 = $ws.Range("B20").Value()
$tc2_D20 = $ws.Range("D20").Value()

$tc3_B28 = $ws.Range("B28").Value()
$tc3_D28 = $ws.Range("D28").Value()

$tc4_B36 = $ws.Range("B36").Value()
$tc4_D36 = $ws.Range("D36").Value()

# 3. Rotate: TC2 <- TC4, TC3 <- TC2(old), TC4 <- TC3(old)
$ws.Range("B20").Value = $tc4_B36
$ws.Range("D20").Value = $tc4_D36

$ws.Range("B28").Value = $tc2_B20
$ws.Range("D28").Value = $tc2_D20

$ws.Range("B36").Value = $tc3_B28
$ws.Range("D36").Value = $tc3_D28
